$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 43 - this pushes the existing rows 43-63
# down to 44-64 (and the sheet dimension grows from A1:R63 to A1:R64).
$ws.Rows(43).Insert()

# Populate the freshly-inserted row 43 with the new weekly price entry.
# Columns that repeat the prior week's descriptive values (market,
# region, category, etc.) are carried over; the date and price columns
# hold the new figures.
$ws.Cells.Item(43, 1).Value = 11
$ws.Cells.Item(43, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(43, 3).Value = "Bíobío"
$ws.Cells.Item(43, 4).Value = 44777
$ws.Cells.Item(43, 5).Value = 8
$ws.Cells.Item(43, 6).Value = 100112037
$ws.Cells.Item(43, 7).Value = "Cebollín"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 150
$ws.Cells.Item(43, 11).Value = 7500
$ws.Cells.Item(43, 12).Value = 8000
$ws.Cells.Item(43, 13).Value = 7767
$ws.Cells.Item(43, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(43, 15).Value = "Región Metropolitana"
$ws.Cells.Item(43, 16).Value = 216
$ws.Cells.Item(43, 17).Value = 36
$ws.Cells.Item(43, 18).Value = "Hortaliza"
